$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. The existing "2022-Q2" sheet (index 2) will become the new
#    "2022-Q3" sheet (new fund-holding detail data). Before overwriting
#    it, duplicate it so its current data survives as a separate sheet
#    named "2022-Q2" (placed right after the new "2022-Q3" sheet).
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ = $wb.Worksheets.Item(2)

$wsQ.Copy($null, $wsQ)
$wsOldCopy = $wb.Worksheets.Item(3)

$wsQ.Name = "2022-Q3"
$wsOldCopy.Name = "2022-Q2"

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: rename the first data row from
#    2022-Q2 to 2022-Q3 with the new totals, and append a new row for
#    the (now second) 2022-Q2 entry with the original totals.
# ------------------------------------------------------------------
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 8
$wsTotal.Range("D2").Value = 0.13

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 0.06
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Replace the contents of the "2022-Q3" sheet with the new fund
#    holdings detail table.
# ------------------------------------------------------------------
$wsQ.Cells.Clear()

# Header row (copy the header style used on the "总计" sheet, s=2)
$wsQ.Range("B1").Value = "基金代码"
$wsQ.Range("C1").Value = "基金名称"
$wsQ.Range("D1").Value = "基金规模"
$wsQ.Range("E1").Value = "股票总仓位"
$wsQ.Range("F1").Value = "仓位占比"
$wsQ.Range("G1").Value = "持有市值(亿元)"
$wsQ.Range("H1").Value = "仓位排名"

$wsTotal.Range("B1:D1").Copy()
$wsQ.Range("B1:D1").PasteSpecial(-4122)
$wsTotal.Range("B1").Copy()
$wsQ.Range("E1").PasteSpecial(-4122)
$wsQ.Range("F1").PasteSpecial(-4122)
$wsQ.Range("G1").PasteSpecial(-4122)
$wsQ.Range("H1").PasteSpecial(-4122)

# Data rows: columns B, D, E, F, G hold numeric-looking text, so force
# a text number format before assignment to keep them as strings
# (e.g. "1.60", "0.0310") instead of auto-converted numbers.
$textRangeB = $wsQ.Range("B2:B9")
$textRangeDG = $wsQ.Range("D2:G9")
$textRangeB.NumberFormat = "@"
$textRangeDG.NumberFormat = "@"

$fund1  = @(0,"003142","鹏华弘达灵活配置混合A","1.60","31.16","1.94","0.0310",7)
$fund2  = @(1,"001331","鹏华弘信灵活配置混合A","1.66","31.16","1.64","0.0272",9)
$fund3  = @(2,"001325","鹏华弘和灵活配置混合A","1.55","28.62","1.38","0.0214",10)
$fund4  = @(3,"001326","鹏华弘和灵活配置混合C","1.43","28.62","1.38","0.0197",10)
$fund5  = @(4,"001327","鹏华弘华灵活配置混合A","0.59","56.25","3.21","0.0189",8)
$fund6  = @(5,"001332","鹏华弘信灵活配置混合C","0.65","31.16","1.64","0.0107",9)
$fund7  = @(6,"003143","鹏华弘达灵活配置混合C","0.13","31.16","1.94","0.0025",7)
$fund8  = @(7,"001328","鹏华弘华灵活配置混合C","0.02","56.25","3.21","0.0006",8)

$funds = @($fund1, $fund2, $fund3, $fund4, $fund5, $fund6, $fund7, $fund8)

$r = 2
foreach ($f in $funds) {
    $wsQ.Cells.Item($r, 1).Value = $f[0]
    $wsQ.Cells.Item($r, 2).Value = $f[1]
    $wsQ.Cells.Item($r, 3).Value = $f[2]
    $wsQ.Cells.Item($r, 4).Value = $f[3]
    $wsQ.Cells.Item($r, 5).Value = $f[4]
    $wsQ.Cells.Item($r, 6).Value = $f[5]
    $wsQ.Cells.Item($r, 7).Value = $f[6]
    $wsQ.Cells.Item($r, 8).Value = $f[7]
    $r = $r + 1
}

# Style the A-column index cells like the "总计" sheet's A2 cell.
$wsTotal.Range("A2").Copy()
$aRange = $wsQ.Range("A2:A9")
$aRange.PasteSpecial(-4122)

# Drop the temporary text-format marker now that values are committed,
# restoring the cells to the default (unstyled) look used in the source.
$textRangeB.ClearFormats()
$textRangeDG.ClearFormats()

Write-Output "done"
